$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.382.76"
$ws.Range("E2").Value = "  +2.17%  "

$ws.Range("D3").Value = "2.998.81"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.18"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.03"
$ws.Range("E6").Value = "  +3.84%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +0.86%  "

$ws.Range("D9").Value = "2.985.18"
$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("E10").Value = "  +2.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.21"
$ws.Range("E11").Value = "  +6.62%  "

$ws.Range("E12").Value = "  +2.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  +2.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.78"
$ws.Range("E14").Value = "  +2.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.122"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.492.93"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.30"
$ws.Range("E17").Value = "  +6.44%  "

$ws.Range("D18").Value = "2.998.62"
$ws.Range("E18").Value = "  +0.92%  "

$ws.Range("D19").Value = "59.355.89"
$ws.Range("E19").Value = "  +2.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.08"
$ws.Range("E20").Value = "  +2.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("E21").Value = "  +2.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.723"
$ws.Range("E22").Value = "  +4.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.14"
$ws.Range("E23").Value = "  +1.78%  "

$ws.Range("E24").Value = "  +3.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.32"
$ws.Range("E25").Value = "  +0.71%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.24"
$ws.Range("E27").Value = "  +10.59%  "

$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.55"
$ws.Range("E29").Value = "  +1.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.91"
$ws.Range("E30").Value = "  +3.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.79"
$ws.Range("E31").Value = "  +1.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.14"
$ws.Range("E32").Value = "  +0.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.101"
$ws.Range("E33").Value = "  +0.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.98"
$ws.Range("E34").Value = "  +5.13%  "

$ws.Range("E35").Value = "  +5.56%  "

$ws.Range("D36").Value = "0.0₃0766"
$ws.Range("E36").Value = "  +8.99%  "

$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.97"
$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("E40").Value = "  +5.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "405.51"
$ws.Range("E41").Value = "  +6.69%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.781.60"
$ws.Range("E42").Value = "  +3.29%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0353"
$ws.Range("E43").Value = "  +0.40%  "

$ws.Range("E44").Value = "  -1.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.253"
$ws.Range("E45").Value = "  +4.18%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.90"
$ws.Range("E47").Value = "  +21.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.53"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("E50").Value = "  +0.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.57"
$ws.Range("E51").Value = "  -0.23%  "
